# REPORTGEN-1102: part 1, added and removed counts missing when no previous snapshot selected
#
# Appends ",EVOLUTION=true" to the RepGen table directives that drive the
# quality-standard evolution table and the per-rule statistics-ratio tables,
# so the report generator also renders added/removed counts when no
# previous snapshot was selected.

$wb = $excel.ActiveWorkbook

# --- Summary sheet: QUALITY_STANDARDS_EVOLUTION table directive ---
$wsSummary = $wb.Worksheets.Item("Summary")
$cell = $wsSummary.Range("B15")
$cell.Value2 = $cell.Value2 + ",EVOLUTION=true"

# --- Each "AN-2013" sheet: RULES_LIST_STATISTICS_RATIO table directive (cell A3) ---
$owaspSheets = @(
    "A1-2013",
    "A2-2013",
    "A3-2013",
    "A4-2013",
    "A5-2013",
    "A6-2013",
    "A7-2013",
    "A8-2013",
    "A9-2013",
    "A10-2013"
)

foreach ($sheetName in $owaspSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range("A3")
    $cell.Value2 = $cell.Value2 + ",EVOLUTION=true"
}
